$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final (post-edit) roster table for rows 2-19 (A: Player, B: Position, C: Team)
$data = @(
    @("Kyrie Irving","PG,SG","Dallas Mavericks"),
    @("CJ McCollum","PG,SG","New Orleans Pelicans"),
    @("Shai Gilgeous-Alexander","PG,SG","Oklahoma City Thunder"),
    @("Zach LaVine","SG,SF","Chicago Bulls"),
    @("RJ Barrett","SG,SF,PF","Toronto Raptors"),
    @("Jaylin Williams","PF,C","Oklahoma City Thunder"),
    @("Jimmy Butler","SF,PF","Miami Heat"),
    @("Jalen Williams","SG,SF,PF,C","Oklahoma City Thunder"),
    @("Nikola Jovic","PF,C","Miami Heat"),
    @("Cason Wallace","PG,SG","Oklahoma City Thunder"),
    @("Isaiah Stewart","PF,C","Detroit Pistons"),
    @("Jordan Poole","PG,SG","Washington Wizards"),
    @("Andrew Nembhard","PG,SG","Indiana Pacers"),
    @("Toumani Camara","SF,PF","Portland Trail Blazers"),
    @("Tobias Harris","SF,PF","Detroit Pistons"),
    @("Lauri Markkanen","SF,PF","Utah Jazz"),
    @("John Collins","PF,C","Utah Jazz"),
    @("Joel Embiid","C","Philadelphia 76ers")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
